$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187295913696289
$ws.Range("B1").Value = 2.35515284538269
$ws.Range("C1").Value = 4.012402057647705
$ws.Range("D1").Value = 2.897143363952637
$ws.Range("E1").Value = 1.131458520889282
